$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Ají" (Inferno, Primera) was inserted before
# the current row 195, shifting the existing rows 195-200 down to 196-201.
$ws.Rows.Item(195).Insert()

# Match the date formatting used by the surrounding rows in column D.
$ws.Range("D195").NumberFormat = $ws.Range("D196").NumberFormat

# Populate the newly inserted row with the new record's data.
$ws.Range("A195").Value = 8
$ws.Range("B195").Value = "Terminal La Palmera de La Serena"
$ws.Range("C195").Value = "Coquimbo"
$ws.Range("D195").Value = 44628
$ws.Range("E195").Value = 4
$ws.Range("F195").Value = 100112021
$ws.Range("G195").Value = "Ají"
$ws.Range("H195").Value = "Inferno"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 500
$ws.Range("K195").Value = 15000
$ws.Range("L195").Value = 16000
$ws.Range("M195").Value = 15500
$ws.Range("N195").Value = "`$/caja 15 kilos"
$ws.Range("O195").Value = "Provincia de Limarí"
$ws.Range("P195").Value = 1033
$ws.Range("Q195").Value = 15
$ws.Range("R195").Value = "Hortaliza"
